$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "558.95") would be
# auto-coerced to a numeric value by plain Range.Value assignment. Force them
# to stay text (matching the source inlineStr cells) by applying a text
# number format for the assignment, then restoring the default cell style so
# no visible style/format change is left behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "61.853.03"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "2.407.03"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "558.95"
$ws.Range("E5").Value = "  +2.52%  "
Set-TextValue $ws.Range("D6") "137.84"
$ws.Range("E6").Value = "  +5.30%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.584"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "2.404.36"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("E10").Value = "  +2.79%  "
Set-TextValue $ws.Range("D11") "5.73"
$ws.Range("E11").Value = "  +4.13%  "
$ws.Range("E12").Value = "  -0.17%  "
Set-TextValue $ws.Range("D13") "0.347"
$ws.Range("E13").Value = "  +3.70%  "
Set-TextValue $ws.Range("D14") "25.77"
$ws.Range("E14").Value = "  +8.91%  "
$ws.Range("D15").Value = "2.839.91"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("D16").Value = "61.850.32"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").Value = "  +4.49%  "
$ws.Range("D18").Value = "2.406.24"
$ws.Range("E18").Value = "  +2.24%  "
Set-TextValue $ws.Range("D19") "11.07"
$ws.Range("E19").Value = "  +4.59%  "
Set-TextValue $ws.Range("D20") "342.60"
$ws.Range("E20").Value = "  +9.27%  "
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("E22").Value = "  +2.99%  "
$ws.Range("E23").Value = "  +0.32%  "
Set-TextValue $ws.Range("D24") "64.91"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  +0.09%  "
Set-TextValue $ws.Range("D27") "8.30"
$ws.Range("E27").Value = "  +5.73%  "
Set-TextValue $ws.Range("D28") "1.50"
$ws.Range("E28").Value = "  +11.74%  "
$ws.Range("E29").Value = "  +15.18%  "
Set-TextValue $ws.Range("D30") "1.79"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").Value = "0.0₃0774"
$ws.Range("E31").Value = "  +5.85%  "
Set-TextValue $ws.Range("D32") "6.33"
$ws.Range("E32").Value = "  +6.61%  "
Set-TextValue $ws.Range("D33") "170.82"
$ws.Range("E33").Value = "  -1.55%  "
Set-TextValue $ws.Range("D34") "0.396"
$ws.Range("E34").Value = "  +3.98%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D35") "1.40"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D36") "374.08"
$ws.Range("E36").Value = "  +15.64%  "
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("E38").Value = "  +11.01%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +7.95%  "
Set-TextValue $ws.Range("D42") "39.05"
$ws.Range("E42").Value = "  +2.95%  "
Set-TextValue $ws.Range("D43") "144.17"
$ws.Range("E43").Value = "  +4.69%  "
$ws.Range("E44").Value = "  +4.52%  "
$ws.Range("E45").Value = "  +7.68%  "
Set-TextValue $ws.Range("D46") "0.0526"
$ws.Range("E46").Value = "  +6.25%  "
Set-TextValue $ws.Range("D47") "0.0956"
$ws.Range("E47").Value = "  +1.71%  "
Set-TextValue $ws.Range("D48") "0.584"
$ws.Range("E48").Value = "  +4.38%  "
Set-TextValue $ws.Range("D49") "17.87"
$ws.Range("E49").Value = "  +5.57%  "
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").Value = "0.0₆0219"
$ws.Range("E51").Value = "  +2.63%  "
